# Rename the header row column captions:
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace '_old$', '_FV2310')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace '_new$', '_FV2404')
}

# Freeze the header row (split below row 1, keep left pane empty so
# the active pane is bottomLeft, matching a one-row freeze).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table ("ListObject") so it gets an
# autofilter and banded-row styling, like the source workbook does.
$rng = $ws.Range("A1:U94")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
